# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the Overview roll-up (zh-cn / de-de status
# columns) and each language sheet's own Status cell get the new
# status text, and the associated "generated at" timestamps advance
# by under a minute (the handoff-report generation run).
#
# Because "Ready for handoff" (17 chars) is wider than "In Translation"
# (14 chars), the Status columns were re-sized when the report was
# regenerated - reflected below as explicit Status column widening on
# all three sheets. Target stored column width is ~17.216 "characters";
# ColumnWidth assignments snap to this host's character-width grid, so
# 16.333333333333332 is the closest input that lands on the nearest
# representable width to the captured value.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 10:43:34"

# Status columns (zh-cn / de-de) widened to fit the new status text.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 10:43:30"
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 10:43:34"
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
